# Fruta / hortaliza, semanal
# Insert 2 new weekly data rows (rows 54-55) into the "Ají" sheet, pushing
# the existing data rows down by two (old row 54 -> new row 56, ...,
# old row 141 -> new row 143). New dimension becomes A1:R143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 54; this shifts every
# row currently at 54..141 down to 56..143 and keeps the date-column style
# (s="2") that Excel propagates on row insert.
$ws.Rows("54:55").Insert()

# --- New row 54 ---------------------------------------------------------
$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44533
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112021
$ws.Range("G54").Value = "Ají"
$ws.Range("H54").Value = "Inferno"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 800
$ws.Range("K54").Value = 11500
$ws.Range("L54").Value = 12000
$ws.Range("M54").Value = 11750
$ws.Range("N54").Value = "`$/caja 12 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 979
$ws.Range("Q54").Value = 12
$ws.Range("R54").Value = "Hortaliza"

# --- New row 55 ---------------------------------------------------------
$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "Terminal La Palmera de La Serena"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44533
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100112021
$ws.Range("G55").Value = "Ají"
$ws.Range("H55").Value = "Inferno"
$ws.Range("I55").Value = "Segunda"
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 7500
$ws.Range("L55").Value = 8000
$ws.Range("M55").Value = 7750
$ws.Range("N55").Value = "`$/caja 12 kilos"
$ws.Range("O55").Value = "Región de Arica y Parinacota"
$ws.Range("P55").Value = 646
$ws.Range("Q55").Value = 12
$ws.Range("R55").Value = "Hortaliza"

Write-Output "Inserted rows 54-55; dimension now A1:R143"
